$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 21749.75
$ws.Range("I113").Value = 21999.5
$ws.Range("J113").Value = 21500
$ws.Range("K113").Value = 21999.5
$ws.Range("L113").Value = 21500
$ws.Range("M113").Value = -18745.5
$ws.Range("N113").Value = -28008
$ws.Range("H116").Value = 11386274
$ws.Range("I116").Value = 15654896
$ws.Range("J116").Value = 3281
$ws.Range("K116").Value = 15654896
$ws.Range("L116").Value = 3281
$ws.Range("M116").Value = -15651454
$ws.Range("N116").Value = -10165
$ws.Range("H127").Value = 666.5
$ws.Range("I127").Value = 666.5
$ws.Range("K127").Value = 1999.5
$ws.Range("M127").Value = 2960.5
$ws.Range("H132").Value = 2503643.8
$ws.Range("I132").Value = 3524.2856
$ws.Range("J132").Value = 20004480
$ws.Range("K132").Value = 10572.8568
$ws.Range("L132").Value = 60013440
$ws.Range("M132").Value = -8042.856800000001
$ws.Range("N132").Value = -60018500
$ws.Range("H133").Value = 99443.5
$ws.Range("J133").Value = 99443.5
$ws.Range("L133").Value = 99443.5
$ws.Range("N133").Value = -109563.5
$ws.Range("H135").Value = 7976.5884
$ws.Range("I135").Value = 9328.385
$ws.Range("K135").Value = 83955.465
$ws.Range("M135").Value = -81420.465
$ws.Range("H138").Value = 308554.53
$ws.Range("I138").Value = 527440.8
$ws.Range("J138").Value = 4545.8335
$ws.Range("K138").Value = 1582322.4
$ws.Range("L138").Value = 13637.5005
$ws.Range("M138").Value = -1577182.4
$ws.Range("N138").Value = -23917.5005
$ws.Range("H141").Value = 6664.9443
$ws.Range("I141").Value = 6498.0625
$ws.Range("K141").Value = 19494.1875
$ws.Range("M141").Value = -14314.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7010.193
$ws.Range("I32").Value = 6878.9243
$ws.Range("J32").Value = 8749.5
$ws.Range("K32").Value = 6878.9243
$ws.Range("L32").Value = 8749.5
$ws.Range("M32").Value = -6591.9243
$ws.Range("N32").Value = -9323.5
$ws.Range("H45").Value = 52459.195
$ws.Range("I45").Value = 83292
$ws.Range("K45").Value = 83292
$ws.Range("M45").Value = -82915
$ws.Range("H61").Value = 9077.074000000001
$ws.Range("I61").Value = 10531.315
$ws.Range("K61").Value = 10531.315
$ws.Range("M61").Value = -10319.315
$ws.Range("H97").Value = 7147002.5
$ws.Range("I97").Value = 6428.8237
$ws.Range("J97").Value = 18182434
$ws.Range("K97").Value = 6428.8237
$ws.Range("L97").Value = 18182434
$ws.Range("M97").Value = -5932.8237
$ws.Range("N97").Value = -18183426
$ws.Range("H102").Value = 9804.308000000001
$ws.Range("I102").Value = 13728.883
$ws.Range("K102").Value = 13728.883
$ws.Range("M102").Value = -12106.883
$ws.Range("H122").Value = 1004563.1
$ws.Range("I122").Value = 3931.1904
$ws.Range("K122").Value = 11793.5712
$ws.Range("M122").Value = -9343.5712
$ws.Range("H132").Value = 3227.6287
$ws.Range("I132").Value = 2965.3076
$ws.Range("J132").Value = 3985.4443
$ws.Range("K132").Value = 8895.9228
$ws.Range("L132").Value = 11956.3329
$ws.Range("M132").Value = -6365.9228
$ws.Range("N132").Value = -17016.3329
$ws.Range("H136").Value = 9077.074000000001
$ws.Range("I136").Value = 10531.315
$ws.Range("K136").Value = 31593.945
$ws.Range("M136").Value = -29043.945

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10033.634
$ws.Range("I94").Value = 15015.277
$ws.Range("J94").Value = 2561.1667
$ws.Range("K94").Value = 15015.277
$ws.Range("L94").Value = 2561.1667
$ws.Range("M94").Value = -14564.277
$ws.Range("N94").Value = -3463.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6696.8823
$ws.Range("J31").Value = 4762.2144
$ws.Range("L31").Value = 4762.2144
$ws.Range("N31").Value = -5352.2144
$ws.Range("H34").Value = 6696.8823
$ws.Range("J34").Value = 4762.2144
$ws.Range("L34").Value = 4762.2144
$ws.Range("N34").Value = -5166.2144
$ws.Range("H86").Value = 10971.286
$ws.Range("I86").Value = 8960.200000000001
$ws.Range("K86").Value = 8960.200000000001
$ws.Range("M86").Value = -7837.200000000001
$ws.Range("H89").Value = 10971.286
$ws.Range("I89").Value = 8960.200000000001
$ws.Range("K89").Value = 44801
$ws.Range("M89").Value = -39185
$ws.Range("H122").Value = 10448.409
$ws.Range("I122").Value = 8737.1875
$ws.Range("K122").Value = 26211.5625
$ws.Range("M122").Value = -23761.5625
$ws.Range("H132").Value = 2798.1538
$ws.Range("I132").Value = 2687.7
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 8063.099999999999
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -5533.099999999999
$ws.Range("N132").Value = -14558.9999
$ws.Range("H133").Value = 80000.5
$ws.Range("J133").Value = 80000.5
$ws.Range("L133").Value = 80000.5
$ws.Range("N133").Value = -85060.5
$ws.Range("H141").Value = 358586.47
$ws.Range("J141").Value = 424939.12
$ws.Range("L141").Value = 424939.12
$ws.Range("N141").Value = -435299.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 556396.9399999999
$ws.Range("I5").Value = 183.36363
$ws.Range("J5").Value = 1430446.9
$ws.Range("K5").Value = 550.0908899999999
$ws.Range("L5").Value = 4291340.699999999
$ws.Range("M5").Value = -438.0908899999999
$ws.Range("N5").Value = -4291564.699999999
$ws.Range("H133").Value = 13991.429
$ws.Range("I133").Value = 5982.6665
$ws.Range("K133").Value = 17947.9995
$ws.Range("M133").Value = -12887.9995
$ws.Range("H135").Value = 556396.9399999999
$ws.Range("I135").Value = 183.36363
$ws.Range("J135").Value = 1430446.9
$ws.Range("K135").Value = 1650.27267
$ws.Range("L135").Value = 12874022.1
$ws.Range("M135").Value = 884.7273299999999
$ws.Range("N135").Value = -12879092.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 10446.615
$ws.Range("I113").Value = 13653.333
$ws.Range("J113").Value = 3231.5
$ws.Range("K113").Value = 13653.333
$ws.Range("L113").Value = 3231.5
$ws.Range("M113").Value = -11483.333
$ws.Range("N113").Value = -7571.5
$ws.Range("H122").Value = 11958.65
$ws.Range("I122").Value = 9418.299999999999
$ws.Range("K122").Value = 28254.9
$ws.Range("M122").Value = -25804.9
$ws.Range("H126").Value = 8942.645500000001
$ws.Range("I126").Value = 10276.134
$ws.Range("J126").Value = 7692.5
$ws.Range("K126").Value = 30828.402
$ws.Range("L126").Value = 23077.5
$ws.Range("M126").Value = -28358.402
$ws.Range("N126").Value = -28017.5
$ws.Range("H132").Value = 4077.0715
$ws.Range("I132").Value = 4187.222
$ws.Range("J132").Value = 3416.1667
$ws.Range("K132").Value = 12561.666
$ws.Range("L132").Value = 10248.5001
$ws.Range("M132").Value = -10031.666
$ws.Range("N132").Value = -15308.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 44265.91
$ws.Range("I7").Value = 87785.39999999999
$ws.Range("K7").Value = 87785.39999999999
$ws.Range("M7").Value = -87673.39999999999
$ws.Range("H40").Value = 22490.822
$ws.Range("I40").Value = 26764.389
$ws.Range("K40").Value = 26764.389
$ws.Range("M40").Value = -26628.389
$ws.Range("H93").Value = 4860.6665
$ws.Range("I93").Value = 6580.5713
$ws.Range("J93").Value = 1420.8572
$ws.Range("K93").Value = 6580.5713
$ws.Range("L93").Value = 1420.8572
$ws.Range("M93").Value = -5332.5713
$ws.Range("N93").Value = -3916.8572
$ws.Range("H122").Value = 4595.722
$ws.Range("I122").Value = 4522.2
$ws.Range("J122").Value = 4762.8184
$ws.Range("K122").Value = 13566.6
$ws.Range("L122").Value = 14288.4552
$ws.Range("M122").Value = -11116.6
$ws.Range("N122").Value = -19188.4552
$ws.Range("H126").Value = 44265.91
$ws.Range("I126").Value = 87785.39999999999
$ws.Range("K126").Value = 263356.2
$ws.Range("M126").Value = -260886.2
$ws.Range("H132").Value = 787318
$ws.Range("I132").Value = 1147984.5
$ws.Range("J132").Value = 5874
$ws.Range("K132").Value = 3443953.5
$ws.Range("L132").Value = 17622
$ws.Range("M132").Value = -3441423.5
$ws.Range("N132").Value = -22682
$ws.Range("H136").Value = 3956.3333
$ws.Range("I136").Value = 2554.6
$ws.Range("J136").Value = 5358.067
$ws.Range("K136").Value = 7663.799999999999
$ws.Range("L136").Value = 16074.201
$ws.Range("M136").Value = -5113.799999999999
$ws.Range("N136").Value = -21174.201

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4765
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4710
$ws.Range("H113").Value = 1727.4783
$ws.Range("J113").Value = 2984.2
$ws.Range("L113").Value = 8952.599999999999
$ws.Range("N113").Value = -13292.6
$ws.Range("H122").Value = 4181.04
$ws.Range("I122").Value = 2236.606
$ws.Range("J122").Value = 7955.5293
$ws.Range("K122").Value = 6709.818000000001
$ws.Range("L122").Value = 23866.5879
$ws.Range("M122").Value = -4259.818000000001
$ws.Range("N122").Value = -28766.5879
$ws.Range("H126").Value = 30244.133
$ws.Range("I126").Value = 46185.668
$ws.Range("J126").Value = 6331.8335
$ws.Range("K126").Value = 138557.004
$ws.Range("L126").Value = 18995.5005
$ws.Range("M126").Value = -136087.004
$ws.Range("N126").Value = -23935.5005
$ws.Range("H132").Value = 19644.61
$ws.Range("I132").Value = 31471.615
$ws.Range("J132").Value = 4269.5
$ws.Range("K132").Value = 94414.845
$ws.Range("L132").Value = 12808.5
$ws.Range("M132").Value = -91884.845
$ws.Range("N132").Value = -17868.5
$ws.Range("H140").Value = 96198.60000000001
$ws.Range("J140").Value = 96198.60000000001
$ws.Range("L140").Value = 96198.60000000001
$ws.Range("N140").Value = -106558.6
